# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.441.30'
$ws.Range("E2").Value = '  -0.95%  '
$ws.Range("D3").Value = '1.827.27'
$ws.Range("E3").Value = '  -1.96%  '
$ws.Range("D4").Value = '''1.007'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.54%  '
$ws.Range("D5").Value = '''331.38'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.47%  '
$ws.Range("D6").Value = '''1.007'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.32%  '
$ws.Range("D7").Value = '''0.4581'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.62%  '
$ws.Range("E8").Value = '  -1.88%  '
$ws.Range("D9").Value = '''46.43'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("E10").Value = '  -0.69%  '
$ws.Range("D11").Value = '''0.9692'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -2.98%  '
$ws.Range("D12").Value = '''21.03'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -2.14%  '
$ws.Range("D13").Value = '1.836.68'
$ws.Range("E13").Value = '  -1.50%  '
$ws.Range("D14").Value = '''5.869'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.95%  '
$ws.Range("D15").Value = '''7.037'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.79%  '
$ws.Range("E16").Value = '  -0.49%  '
$ws.Range("D17").Value = '''88.17'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.26%  '
$ws.Range("D18").Value = '''0.06655'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.65%  '
$ws.Range("E19").Value = '  -1.22%  '
$ws.Range("D20").Value = '''17.15'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.46%  '
$ws.Range("D21").Value = '''1.007'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.37%  '
$ws.Range("D22").Value = '27.438.30'
$ws.Range("E22").Value = '  -0.92%  '
$ws.Range("D23").Value = '''5.315'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.62%  '
$ws.Range("E24").Value = '  -0.95%  '
$ws.Range("D25").Value = '''2.302'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.87%  '
$ws.Range("B26").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C26").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D26").Value = '2.065.29'
$ws.Range("E26").Value = '  -1.02%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '''157.09'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.53%  '
$ws.Range("D28").Value = '''19.36'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.58%  '
$ws.Range("D29").Value = '''2.060'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.96%  '
$ws.Range("D30").Value = '''5.224'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.87%  '
$ws.Range("D31").Value = '''118.29'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.26%  '
$ws.Range("D32").Value = '''0.9444'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.56%  '
$ws.Range("D33").Value = '''0.09296'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.50%  '
$ws.Range("D34").Value = '''3.597'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.09%  '
$ws.Range("D35").Value = '''5.236'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.10%  '
$ws.Range("D36").Value = '''1.310'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.82%  '
$ws.Range("D37").Value = '''0.05930'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.62%  '
$ws.Range("D38").Value = '''0.02185'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.04%  '
$ws.Range("D39").Value = '''1.159'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -3.36%  '
$ws.Range("D40").Value = '''7.971'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.32%  '
$ws.Range("D41").Value = '''0.5782'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.14%  '
$ws.Range("D42").Value = '''0.1834'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.38%  '
$ws.Range("D43").Value = '''10.00'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.81%  '
$ws.Range("D44").Value = '''1.274'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.51%  '
$ws.Range("D45").Value = '''0.5482'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.29%  '
$ws.Range("D46").Value = '''11.94'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.82%  '
$ws.Range("D47").Value = '''1.865'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.50%  '
$ws.Range("D48").Value = '''0.06643'
$ws.Range("D48").ClearFormats()
$ws.Range("D49").Value = '''109.84'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.16%  '
$ws.Range("D50").Value = '''1.038'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.07%  '
$ws.Range("E51").Value = '  -0.32%  '
